$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 835.9
$ws.Range("I9").Value = 839.875
$ws.Range("K9").Value = 839.875
$ws.Range("M9").Value = -670.875

$ws.Range("H11").Value = 72.92308
$ws.Range("I11").Value = 72.92308
$ws.Range("K11").Value = 72.92308
$ws.Range("M11").Value = 67.07692

$ws.Range("H33").Value = 422.90475
$ws.Range("I33").Value = 431.94736
$ws.Range("J33").Value = 337
$ws.Range("K33").Value = 431.94736
$ws.Range("L33").Value = 337
$ws.Range("M33").Value = -202.94736
$ws.Range("N33").Value = -795

$ws.Range("H40").Value = 3954.6667
$ws.Range("I40").Value = 1749.75
$ws.Range("J40").Value = 5718.6
$ws.Range("K40").Value = 1749.75
$ws.Range("L40").Value = 5718.6
$ws.Range("M40").Value = -1574.75
$ws.Range("N40").Value = -6068.6

$ws.Range("H62").Value = 3883
$ws.Range("I62").Value = 3766
$ws.Range("K62").Value = 3766
$ws.Range("M62").Value = -3142

$ws.Range("H65").Value = 3883
$ws.Range("I65").Value = 3766
$ws.Range("K65").Value = 18830
$ws.Range("M65").Value = -15710

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2777.5
$ws.Range("I2").Value = 2777.5
$ws.Range("K2").Value = 2777.5
$ws.Range("M2").Value = -2664.5

$ws.Range("H19").Value = 3750
$ws.Range("I19").Value = 3750
$ws.Range("K19").Value = 3750
$ws.Range("M19").Value = -3521

$ws.Range("H97").Value = 1536.2354
$ws.Range("I97").Value = 1055.3077
$ws.Range("J97").Value = 3099.25
$ws.Range("K97").Value = 1055.3077
$ws.Range("L97").Value = 3099.25
$ws.Range("M97").Value = -559.3077000000001
$ws.Range("N97").Value = -4091.25

$ws.Range("H102").Value = 1507
$ws.Range("I102").Value = 1507
$ws.Range("K102").Value = 1507
$ws.Range("M102").Value = 115

$ws.Range("H116").Value = 2777.5
$ws.Range("I116").Value = 2777.5
$ws.Range("K116").Value = 2777.5
$ws.Range("M116").Value = -483.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2777.5
$ws.Range("I3").Value = 2777.5
$ws.Range("K3").Value = 2777.5
$ws.Range("M3").Value = -2663.5

$ws.Range("H11").Value = 1624.25
$ws.Range("J11").Value = 1665.6666
$ws.Range("L11").Value = 1665.6666
$ws.Range("N11").Value = -1945.6666

$ws.Range("H19").Value = 18200
$ws.Range("J19").Value = 25000
$ws.Range("L19").Value = 25000
$ws.Range("N19").Value = -25346

$ws.Range("H105").Value = 1957.4
$ws.Range("I105").Value = 2222
$ws.Range("K105").Value = 2222
$ws.Range("M105").Value = -475

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 656.7143
$ws.Range("J15").Value = 999
$ws.Range("L15").Value = 999
$ws.Range("N15").Value = -1339

$ws.Range("H19").Value = 863.36365
$ws.Range("I19").Value = 510.77777
$ws.Range("J19").Value = 2450
$ws.Range("K19").Value = 510.77777
$ws.Range("L19").Value = 2450
$ws.Range("M19").Value = -340.77777
$ws.Range("N19").Value = -2790

$ws.Range("H24").Value = 863.36365
$ws.Range("I24").Value = 510.77777
$ws.Range("J24").Value = 2450
$ws.Range("K24").Value = 510.77777
$ws.Range("L24").Value = 2450
$ws.Range("M24").Value = -340.77777
$ws.Range("N24").Value = -2790

$ws.Range("H32").Value = 6852.5
$ws.Range("I32").Value = 6852.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 6852.5
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H58").Value = 3554.4285
$ws.Range("I58").Value = 2480.1667
$ws.Range("K58").Value = 2480.1667
$ws.Range("M58").Value = -2277.1667

$ws.Range("H86").Value = 3893.5
$ws.Range("I86").Value = 3894.375
$ws.Range("K86").Value = 3894.375
$ws.Range("M86").Value = -2771.375

$ws.Range("H89").Value = 3893.5
$ws.Range("I89").Value = 3894.375
$ws.Range("K89").Value = 19471.875
$ws.Range("M89").Value = -13855.875

$ws.Range("H99").Value = 6162.4
$ws.Range("I99").Value = 6703
$ws.Range("K99").Value = 6703
$ws.Range("M99").Value = -5205

$ws.Range("H126").Value = 6162.4
$ws.Range("I126").Value = 6703
$ws.Range("K126").Value = 20109
$ws.Range("M126").Value = -17639

$ws.Range("H132").Value = 4186.5
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 3356
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 10068
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -15128

$ws.Range("H135").Value = 100000
$ws.Range("I135").Value = 100000
$ws.Range("K135").Value = 100000
$ws.Range("M135").Value = -94930

$ws.Range("H136").Value = 3554.4285
$ws.Range("I136").Value = 2480.1667
$ws.Range("K136").Value = 7440.500100000001
$ws.Range("M136").Value = -4890.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 148
$ws.Range("I8").Value = 148
$ws.Range("K8").Value = 444
$ws.Range("M8").Value = -305

$ws.Range("H18").Value = 501.14285
$ws.Range("I18").Value = 459.66666
$ws.Range("K18").Value = 1378.99998
$ws.Range("M18").Value = -1209.99998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 475
$ws.Range("I2").Value = 480
$ws.Range("J2").Value = 470
$ws.Range("K2").Value = 480
$ws.Range("L2").Value = 470
$ws.Range("M2").Value = -367
$ws.Range("N2").Value = -696

$ws.Range("H3").Value = 25374.143
$ws.Range("I3").Value = 3383
$ws.Range("J3").Value = 41867.5
$ws.Range("K3").Value = 3383
$ws.Range("L3").Value = 41867.5
$ws.Range("M3").Value = -3267
$ws.Range("N3").Value = -42099.5

$ws.Range("H70").Value = 3974.9092
$ws.Range("J70").Value = 3605.5557
$ws.Range("L70").Value = 3605.5557
$ws.Range("N70").Value = -4145.5557

$ws.Range("H73").Value = 3974.9092
$ws.Range("J73").Value = 3605.5557
$ws.Range("L73").Value = 3605.5557
$ws.Range("N73").Value = -5477.5557

$ws.Range("H122").Value = 4887.8887
$ws.Range("I122").Value = 4253.25
$ws.Range("K122").Value = 12759.75
$ws.Range("M122").Value = -10309.75

$ws.Range("H126").Value = 3497.5
$ws.Range("I126").Value = 3583.2856
$ws.Range("J126").Value = 2897
$ws.Range("K126").Value = 10749.8568
$ws.Range("L126").Value = 8691
$ws.Range("M126").Value = -8279.856800000001
$ws.Range("N126").Value = -13631

$ws.Range("H132").Value = 3844.5
$ws.Range("I132").Value = 3821.3333
$ws.Range("J132").Value = 3914
$ws.Range("K132").Value = 11463.9999
$ws.Range("L132").Value = 11742
$ws.Range("M132").Value = -8933.999899999999
$ws.Range("N132").Value = -16802

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2115.5
$ws.Range("I9").Value = 638.6
$ws.Range("K9").Value = 638.6
$ws.Range("M9").Value = -414.6

$ws.Range("H31").Value = 9339.1
$ws.Range("I31").Value = 2979
$ws.Range("K31").Value = 2979
$ws.Range("M31").Value = -2731

$ws.Range("H132").Value = 2952.7334
$ws.Range("I132").Value = 2913.6428
$ws.Range("K132").Value = 8740.928400000001
$ws.Range("M132").Value = -6210.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 475000
$ws.Range("I61").Value = 475000
$ws.Range("K61").Value = 475000
$ws.Range("M61").Value = -474708

$ws.Range("H96").Value = 5245.5557
$ws.Range("I96").Value = 5043
$ws.Range("K96").Value = 5043
$ws.Range("M96").Value = -3670
